$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Range("A2").Value = "CE 544 - Brigham Young University"
$ws.Range("A3").Select() | Out-Null
